# Applies the edits described by the diff:
#  1. Final!A9 team name changes from "ElonGPT o1-mini" to "ElonGPT 4o"
#     (the now-unused "ElonGPT o1-mini" shared string disappears on save).
#  2. Playoffs!B2:D13 projected-score values are updated (B,C,D mirror each row).
#  3. The active sheet / selection moves from Final!D14 to Playoffs!H12.

$wb = $excel.ActiveWorkbook

# --- 1. Update the team name on the Final standings sheet ---------------
$final = $wb.Worksheets.Item("Final")
$final.Range("A9").Value = "ElonGPT 4o"

# --- 2. Update the playoff projection numbers ----------------------------
$playoffs = $wb.Worksheets.Item("Playoffs")

$newValues = @(
    143.19999999999999,
    148.19999999999999,
    134,
    135.80000000000001,
    124.8,
    126.3,
    154.19999999999999,
    128.4,
    147,
    127.2,
    130.4,
    128.6
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $val = $newValues[$i]
    $playoffs.Range("B$row").Value = $val
    $playoffs.Range("C$row").Value = $val
    $playoffs.Range("D$row").Value = $val
}

# --- 3. Move the active sheet / selection to match the saved view state --
$playoffs.Activate()
$playoffs.Range("H12").Select()
